# Adds a new "Exceptional Items" column to the "Quarterly" sheet between the
# existing "P/l before exceptional items & tax" and "P/l before tax" columns.
#
# The source data's row 1 (lower-case header row) is left untouched except
# that a duplicate of the last cell (T1, "Diluted eps.") is appended at U1 -
# this mirrors the shape of the authoritative edit exactly (row 1 itself is
# NOT shifted). Row 2 (the proper header row) gets a brand-new "Exceptional
# Items" label inserted at column L, with every following header shifting
# one column to the right (through U2). Data rows 3-12 get a blank cell at
# the new column L, with their numbers shifted right starting at the old
# column L (so old L -> new M, old M -> new N, ... old T -> new U).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

$lastDataRow = 12
$firstShiftCol = 12   # column L
$lastShiftCol = 20    # column T (last populated column before the edit)

# --- Row 1: just append a duplicate of T1 at U1 (value + style), nothing else moves.
$ws.Cells.Item(1, $lastShiftCol).Copy($ws.Cells.Item(1, $lastShiftCol + 1))

# --- Rows 2-12: shift columns L..T into M..U, working from the right so we
# never clobber a value before it has been copied onward.
for ($r = 2; $r -le $lastDataRow; $r++) {
    for ($c = $lastShiftCol; $c -ge $firstShiftCol; $c--) {
        $ws.Cells.Item($r, $c + 1).Value2 = $ws.Cells.Item($r, $c).Value2
    }
}

# --- New column L: header label on row 2, blank on the data rows.
$ws.Cells.Item(2, $firstShiftCol).Value2 = "Exceptional Items"
for ($r = 3; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, $firstShiftCol).Value2 = $null
}
